$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-19 (A:T) reflecting the refreshed TPM figures
# and the new sending-cluster "ECs" block (rows 14-19).
$data = @(
  @("ECs", "Slitrk1", "Ptprs", "ECs", 1, 0.3333333333333333, 0.006001333333333333, 0.018004, 0.06719314483623444, 0.06719314483623443, 3, 1, 4.260872666666667, 12.782618, 0.09064705929364961, 0.09064705929364959, 0.02557091716355556, 0.230138254472, 0.006090860984096929, 0.006090860984096927),
  @("ECs", "Slitrk1", "Ptprs", "FAPs", 1, 0.3333333333333333, 0.006001333333333333, 0.018004, 0.06719314483623444, 0.06719314483623443, 3, 1, 20.524797, 61.574391, 0.4366505728284585, 0.4366505728284584, 0.123176148396, 1.108585335564, 0.02933992518288735, 0.02933992518288733),
  @("ECs", "Slitrk1", "Ptprs", "Inflammatory-Mac", 1, 0.3333333333333333, 0.006001333333333333, 0.018004, 0.06719314483623444, 0.06719314483623443, 3, 1, 8.931090666666666, 26.793272, 0.190002651698962, 0.1900026516989619, 0.05359845212088889, 0.482386069088, 0.01276687569487696, 0.01276687569487695),
  @("ECs", "Slitrk1", "Ptprs", "MuSCs", 1, 0.3333333333333333, 0.006001333333333333, 0.018004, 0.06719314483623444, 0.06719314483623443, 3, 1, 6.457974333333333, 19.373923, 0.1373888468646722, 0.1373888468646721, 0.03875645663244445, 0.348808109692, 0.009231588686261151, 0.009231588686261146),
  @("ECs", "Slitrk1", "Ptprs", "Neutrophils", 1, 0.3333333333333333, 0.006001333333333333, 0.018004, 0.06719314483623444, 0.06719314483623443, 3, 1, 1.948535, 5.845605, 0.04145370713904261, 0.0414537071390426, 0.01169380804666667, 0.10524427242, 0.002785404947792536, 0.002785404947792535),
  @("ECs", "Slitrk1", "Ptprs", "Resolving-Mac", 1, 0.3333333333333333, 0.006001333333333333, 0.018004, 0.06719314483623444, 0.06719314483623443, 3, 1, 4.881814666666666, 14.645444, 0.1038571621752152, 0.1038571621752152, 0.02929739708622222, 0.263676573776, 0.006978489340319523, 0.006978489340319519),
  @("FAPs", "Slitrk1", "Ptprs", "ECs", 1, 0.3333333333333333, 0.03898133333333333, 0.116944, 0.4364494073388469, 0.4364494073388469, 3, 1, 4.260872666666667, 12.782618, 0.09064705929364961, 0.09064705929364959, 0.1660944977102222, 1.494850479392, 0.03956285530572268, 0.03956285530572268),
  @("FAPs", "Slitrk1", "Ptprs", "FAPs", 1, 0.3333333333333333, 0.03898133333333333, 0.116944, 0.4364494073388469, 0.4364494073388469, 3, 1, 20.524797, 61.574391, 0.4366505728284585, 0.4366505728284584, 0.800083953456, 7.200755581104, 0.1905758837251487, 0.1905758837251487),
  @("FAPs", "Slitrk1", "Ptprs", "Inflammatory-Mac", 1, 0.3333333333333333, 0.03898133333333333, 0.116944, 0.4364494073388469, 0.4364494073388469, 3, 1, 8.931090666666666, 26.793272, 0.190002651698962, 0.1900026516989619, 0.3481458223075555, 3.133312400768, 0.0829265447268213, 0.08292654472682127),
  @("FAPs", "Slitrk1", "Ptprs", "MuSCs", 1, 0.3333333333333333, 0.03898133333333333, 0.116944, 0.4364494073388469, 0.4364494073388469, 3, 1, 6.457974333333333, 19.373923, 0.1373888468646722, 0.1373888468646721, 0.2517404501457778, 2.265664051312, 0.05996328078905377, 0.05996328078905374),
  @("FAPs", "Slitrk1", "Ptprs", "Neutrophils", 1, 0.3333333333333333, 0.03898133333333333, 0.116944, 0.4364494073388469, 0.4364494073388469, 3, 1, 1.948535, 5.845605, 0.04145370713904261, 0.0414537071390426, 0.07595649234666667, 0.68360843112, 0.01809244591283327, 0.01809244591283327),
  @("FAPs", "Slitrk1", "Ptprs", "Resolving-Mac", 1, 0.3333333333333333, 0.03898133333333333, 0.116944, 0.4364494073388469, 0.4364494073388469, 3, 1, 4.881814666666666, 14.645444, 0.1038571621752152, 0.1038571621752152, 0.1902996447928889, 1.712696803136, 0.04532839687926718, 0.04532839687926716),
  @("MuSCs", "Slitrk1", "Ptprs", "ECs", 3, 1, 0.044332, 0.132996, 0.4963574478249187, 0.4963574478249186, 3, 1, 4.260872666666667, 12.782618, 0.09064705929364961, 0.09064705929364959, 0.1888930070586667, 1.700037063528, 0.04499334300383, 0.04499334300382998),
  @("MuSCs", "Slitrk1", "Ptprs", "FAPs", 3, 1, 0.044332, 0.132996, 0.4963574478249187, 0.4963574478249186, 3, 1, 20.524797, 61.574391, 0.4366505728284585, 0.4366505728284584, 0.9099053006040001, 8.189147705436, 0.2167347639204225, 0.2167347639204224),
  @("MuSCs", "Slitrk1", "Ptprs", "Inflammatory-Mac", 3, 1, 0.044332, 0.132996, 0.4963574478249187, 0.4963574478249186, 3, 1, 8.931090666666666, 26.793272, 0.190002651698962, 0.1900026516989619, 0.3959331114346667, 3.563398002912, 0.09430923127726372, 0.09430923127726368),
  @("MuSCs", "Slitrk1", "Ptprs", "MuSCs", 3, 1, 0.044332, 0.132996, 0.4963574478249187, 0.4963574478249186, 3, 1, 6.457974333333333, 19.373923, 0.1373888468646722, 0.1373888468646721, 0.2862949181453334, 2.576654263308, 0.06819397738935726, 0.06819397738935723),
  @("MuSCs", "Slitrk1", "Ptprs", "Neutrophils", 3, 1, 0.044332, 0.132996, 0.4963574478249187, 0.4963574478249186, 3, 1, 1.948535, 5.845605, 0.04145370713904261, 0.0414537071390426, 0.08638245362000001, 0.77744208258, 0.0205758562784168, 0.0205758562784168),
  @("MuSCs", "Slitrk1", "Ptprs", "Resolving-Mac", 3, 1, 0.044332, 0.132996, 0.4963574478249187, 0.4963574478249186, 3, 1, 4.881814666666666, 14.645444, 0.1038571621752152, 0.1038571621752152, 0.2164206078026667, 1.947785470224, 0.05155027595562849, 0.05155027595562846)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
